# Insert a brand-new record as row 69, pushing all existing rows
# (old 69..117) down by one (to 70..118).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(69).Insert()

$ws.Cells.Item(69, 1).Value = 4
$ws.Cells.Item(69, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(69, 3).Value = "Los Lagos"
$ws.Cells.Item(69, 4).Value = 44981
$ws.Cells.Item(69, 5).Value = 10
$ws.Cells.Item(69, 6).Value = "Fruta"
$ws.Cells.Item(69, 7).Value = 100104
$ws.Cells.Item(69, 8).Value = "Frutos de pepita"
$ws.Cells.Item(69, 9).Value = 100104003
$ws.Cells.Item(69, 10).Value = "Membrillo"
$ws.Cells.Item(69, 11).Value = "Champion"
$ws.Cells.Item(69, 12).Value = "Primera"
$ws.Cells.Item(69, 13).Value = 200
$ws.Cells.Item(69, 14).Value = 16000
$ws.Cells.Item(69, 15).Value = 17000
$ws.Cells.Item(69, 16).Value = 16500
$ws.Cells.Item(69, 17).Value = "`$/caja 18 kilos empedrada"
$ws.Cells.Item(69, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(69, 19).Value = 917
$ws.Cells.Item(69, 20).Value = 18
